# Insert a new record (weekly/daily price log) at row 105, pushing all
# subsequent rows down by one (dimension grows from A1:R197 to A1:R198).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(105).Insert()

$ws.Cells.Item(105, 1).Value = 11
$ws.Cells.Item(105, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(105, 3).Value = "Bíobío"
$ws.Cells.Item(105, 4).Value = 44589
$ws.Cells.Item(105, 5).Value = 8
$ws.Cells.Item(105, 6).Value = 100112017
$ws.Cells.Item(105, 7).Value = "Apio"
$ws.Cells.Item(105, 8).Value = "Americana (o)"
$ws.Cells.Item(105, 9).Value = "Primera"
$ws.Cells.Item(105, 10).Value = 220
$ws.Cells.Item(105, 11).Value = 7000
$ws.Cells.Item(105, 12).Value = 7500
$ws.Cells.Item(105, 13).Value = 7273
$ws.Cells.Item(105, 14).Value = "$/docena de matas"
$ws.Cells.Item(105, 15).Value = "Región de Coquimbo"
$ws.Cells.Item(105, 16).Value = 1212
$ws.Cells.Item(105, 17).Value = 6
$ws.Cells.Item(105, 18).Value = "Hortaliza"
